$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "92.844.24"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.108.96"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.48"
$ws.Range("E5").Value = "  -3.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "613.40"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("E7").Value = "  -1.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.391"
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.107.12"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.790"
$ws.Range("E11").Value = "  +5.35%  "
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.570.79"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.09"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.694.01"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.110.04"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.80"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.44"
$ws.Range("E20").Value = "  -2.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.84"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "438.89"
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.09"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("E25").Value = "  +4.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.58"
$ws.Range("E26").Value = "  -4.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "85.54"
$ws.Range("E27").Value = "  -4.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.88"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.276.66"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.182"
$ws.Range("E31").Value = "  +8.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.126"
$ws.Range("E32").Value = "  -9.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.233"
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.14"
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -6.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.10"
$ws.Range("E36").Value = "  +6.25%  "
$ws.Range("E37").Value = "  -7.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.66"
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.98"
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("E40").Value = "  -19.19%  "
$ws.Range("E41").Value = "  +7.83%  "
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "467.84"
$ws.Range("E43").Value = "  -4.36%  "
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "159.19"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.681"
$ws.Range("E48").Value = "  -2.71%  "
$ws.Range("E49").Value = "  -4.42%  "
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.32"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.01"
$ws.Range("E51").Value = "  -0.42%  "
